$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select rows 8:9 the way the author appears to have (A8:XFD9), with
# active cell A8, mirroring the recorded selection in the saved file.
$ws.Range("A8:XFD9").Select()

# Row 8 - Magnesium chloride unit price: replace baseline and the
# formula-driven Lower/Upper bounds with new literal values.
$ws.Range("E8").Value = 0.38
$ws.Range("G8").Value = 0.349
$ws.Range("I8").Value = 0.411

# Row 9 - Zinc sulfate unit price: same treatment.
$ws.Range("E9").Value = 0.795
$ws.Range("G9").Value = 0.657
$ws.Range("I9").Value = 0.931
